$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = 45984
$ws.Range("B35").Value = 78
$ws.Range("C35").Value = 88
$ws.Range("D35").Value = 84

$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)
